# The "reviews_count" column (E) is entirely empty and is being removed
# from the sheet, shifting reviews_average/latitude/longitude/
# is_permanently_closed/gmaps_link/latest_review_date one column to the
# left (F->E, G->F, H->G, I->H, J->I, K->J) and updating the used range
# from A1:K21 to A1:J21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").EntireColumn.Delete()
